# Adds one new weekly price-report block (Melón, "Extra" / "Primera",
# week of 2021-12-27) as two new rows right before the existing data
# block, pushing all the existing weekly blocks down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the top of the data block (row 27 is the
# first data row after the header); everything below shifts down by 2,
# which is exactly what the target sheet shows (old row N -> new row N+2,
# dimension grows from R54 to R56).
$ws.Rows("27:28").Insert()

# New row 27: Melón, Calidad "Extra", week of 2021-12-27 (serial 44557)
$ws.Cells.Item(27,1).Value  = 8
$ws.Cells.Item(27,2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(27,3).Value  = "Coquimbo"
$ws.Cells.Item(27,4).Value  = 44557
$ws.Cells.Item(27,5).Value  = 4
$ws.Cells.Item(27,6).Value  = 100112027
$ws.Cells.Item(27,7).Value  = "Melón"
$ws.Cells.Item(27,8).Value  = "Tuna"
$ws.Cells.Item(27,9).Value  = "Extra"
$ws.Cells.Item(27,10).Value = 5000
$ws.Cells.Item(27,11).Value = 1200
$ws.Cells.Item(27,12).Value = 1250
$ws.Cells.Item(27,13).Value = 1225
$ws.Cells.Item(27,14).Value = "`$/unidad"
$ws.Cells.Item(27,15).Value = "Región de O'Higgins"
$ws.Cells.Item(27,16).Value = 1225
$ws.Cells.Item(27,17).Value = 1
$ws.Cells.Item(27,18).Value = "Hortaliza"

# New row 28: Melón, Calidad "Primera", same week (serial 44557)
$ws.Cells.Item(28,1).Value  = 8
$ws.Cells.Item(28,2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(28,3).Value  = "Coquimbo"
$ws.Cells.Item(28,4).Value  = 44557
$ws.Cells.Item(28,5).Value  = 4
$ws.Cells.Item(28,6).Value  = 100112027
$ws.Cells.Item(28,7).Value  = "Melón"
$ws.Cells.Item(28,8).Value  = "Tuna"
$ws.Cells.Item(28,9).Value  = "Primera"
$ws.Cells.Item(28,10).Value = 4000
$ws.Cells.Item(28,11).Value = 900
$ws.Cells.Item(28,12).Value = 1000
$ws.Cells.Item(28,13).Value = 950
$ws.Cells.Item(28,14).Value = "`$/unidad"
$ws.Cells.Item(28,15).Value = "Región de O'Higgins"
$ws.Cells.Item(28,16).Value = 950
$ws.Cells.Item(28,17).Value = 1
$ws.Cells.Item(28,18).Value = "Hortaliza"
